# Update the cryptocurrency price/volume table to the latest snapshot.
# Most cells in columns D (Price) and E (Volume(1h)) are stored as plain
# text in the workbook (they include thousands separators written with
# dots, percent signs and padding spaces), so we force every write to stay
# text-typed (Excel would otherwise silently coerce clean-looking numeric
# strings like "545.92" into binary doubles, e.g. 545.91999999999996).
# A leading apostrophe forces text entry the same way a user typing into
# the formula bar would, and resetting the Style back to "Normal"
# afterwards clears the implicit @ text-number-format that the apostrophe
# trick applies, so the cell's style index is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

$changes = @(
  @('D2', '60.303.99'),
  @('E2', '  +4.15%  '),
  @('D3', '2.340.84'),
  @('E4', '  +0.00%  '),
  @('D5', '545.92'),
  @('E5', '  +3.03%  '),
  @('D6', '131.90'),
  @('E6', '  +1.04%  '),
  @('E7', '  -0.01%  '),
  @('D8', '0.585'),
  @('E8', '  +0.63%  '),
  @('D9', '2.337.37'),
  @('E9', '  +2.46%  '),
  @('E10', '  +2.01%  '),
  @('D11', '5.52'),
  @('E11', '  +1.18%  '),
  @('E12', '  +0.87%  '),
  @('E13', '  +1.72%  '),
  @('D14', '23.82'),
  @('E14', '  +2.12%  '),
  @('D15', '2.757.30'),
  @('E15', '  +2.50%  '),
  @('D16', '60.221.42'),
  @('E16', '  +4.16%  '),
  @('D18', '2.338.89'),
  @('E18', '  +2.37%  '),
  @('D19', '10.61'),
  @('E19', '  +1.31%  '),
  @('D20', '4.15'),
  @('E20', '  -0.02%  '),
  @('D21', '6.80'),
  @('E21', '  +6.59%  '),
  @('D22', '313.72'),
  @('E22', '  +1.16%  '),
  @('E23', '  -0.23%  '),
  @('D24', '63.50'),
  @('E24', '  +2.07%  '),
  @('D25', '0.171'),
  @('E25', '  +3.09%  '),
  @('E26', '  +0.03%  '),
  @('D27', '7.91'),
  @('E27', '  -0.56%  '),
  @('E28', '  +8.49%  '),
  @('E29', '  +2.97%  '),
  @('D30', '171.72'),
  @('E30', '  +1.20%  '),
  @('E31', '  +13.46%  '),
  @('E32', '  +2.46%  '),
  @('D33', '5.95'),
  @('E33', '  +4.45%  '),
  @('E34', '  +13.77%  '),
  @('D35', '0.383'),
  @('E35', '  +1.91%  '),
  @('D36', '18.01'),
  @('E36', '  +1.70%  '),
  @('E37', '  +0.03%  '),
  @('E38', '  +0.03%  '),
  @('D40', '321.56'),
  @('E40', '  +12.53%  '),
  @('D41', '38.13'),
  @('E41', '  -0.84%  '),
  @('E42', '  +3.17%  '),
  @('D43', '141.20'),
  @('E43', '  +0.46%  '),
  @('E44', '  +1.83%  '),
  @('B45', 'Stellar'),
  @('C45', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'),
  @('D45', '0.0945'),
  @('E45', '  -0.02%  '),
  @('B46', 'InjectiveProtocol'),
  @('C46', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
  @('D46', '19.52'),
  @('E46', '  +9.17%  '),
  @('E47', '  +0.82%  '),
  @('E48', '  +1.71%  '),
  @('E49', '  +2.27%  '),
  @('E50', '  +0.83%  '),
  @('E51', '  +17.65%  ')
)

foreach ($item in $changes) {
    Set-TextValue $ws.Range($item[0]) $item[1]
}
